# Weekly update: insert a new price record for "Haba" as row 31,
# pushing the existing rows 31:59 down to 32:60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a new row at position 31.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with the new weekly record.
$ws.Range("A31").Value = 2
$ws.Range("B31").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = 44763
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = 100112026
$ws.Range("G31").Value = "Haba"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 1100
$ws.Range("K31").Value = 11000
$ws.Range("L31").Value = 12000
$ws.Range("M31").Value = 11500
$ws.Range("N31").Value = "`$/saco 25 kilos"
$ws.Range("O31").Value = "Provincia de Limarí"
$ws.Range("P31").Value = 460
$ws.Range("Q31").Value = 25
$ws.Range("R31").Value = "Hortaliza"
